# Add a new health-bar UI test case row to the Functionality Testing doc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 7

# Text of the new test case (becomes a new shared string entry)
$ws.Cells.Item($newRow, 2).Value = "health bar image fill changes in correspondance with the player's current health value when damge is received"

# Match the thin-border / wrap-text formatting already used elsewhere in the sheet (e.g. B3)
$ws.Cells.Item($newRow, 2).WrapText = $true
$ws.Cells.Item($newRow, 2).Borders.Item(7).LineStyle = 1
$ws.Cells.Item($newRow, 2).Borders.Item(8).LineStyle = 1
$ws.Cells.Item($newRow, 2).Borders.Item(9).LineStyle = 1
$ws.Cells.Item($newRow, 2).Borders.Item(10).LineStyle = 1

# Taller row to accommodate the wrapped, longer description text
$ws.Rows.Item($newRow).RowHeight = 120

# Refresh selection as left by the author after making the edit
$ws.Range("B1:D1048576").Select()
